$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed cryptocurrency price/volume(1h) snapshot values (GitHub Actions
# scraper update). Source data stores every figure as literal text (the
# sheet uses inlineStr cells even for numeric-looking Price/Volume(1h)
# columns), so each cell is forced to Text format before the new value is
# written -- otherwise Excel would auto-coerce a string like "258.46" or
# "0.71%" into a real Number/Percentage and silently change its type.
$updates = @(
    @{ Cell = "D2"; Value = "258.46" }
    @{ Cell = "E2"; Value = "0.71%" }
    @{ Cell = "D3"; Value = "26.88" }
    @{ Cell = "E3"; Value = "-1.51%" }
    @{ Cell = "D4"; Value = "4.644" }
    @{ Cell = "E4"; Value = "1.85%" }
    @{ Cell = "D5"; Value = "0.05981" }
    @{ Cell = "E5"; Value = "1.42%" }
    @{ Cell = "E6"; Value = "0.28%" }
    @{ Cell = "D7"; Value = "0.8559" }
    @{ Cell = "E7"; Value = "-0.37%" }
    @{ Cell = "D8"; Value = "0.9197" }
    @{ Cell = "E8"; Value = "-1.37%" }
    @{ Cell = "D9"; Value = "0.1388" }
    @{ Cell = "E9"; Value = "-1.37%" }
    @{ Cell = "D10"; Value = "0.04704" }
    @{ Cell = "E10"; Value = "29.68%" }
    @{ Cell = "D11"; Value = "0.07017" }
    @{ Cell = "E11"; Value = "-0.94%" }
    @{ Cell = "D12"; Value = "0.03053" }
    @{ Cell = "E12"; Value = "-5.46%" }
    @{ Cell = "D13"; Value = "0.09112" }
    @{ Cell = "E13"; Value = "-1.05%" }
    @{ Cell = "D14"; Value = "0.001527" }
    @{ Cell = "E14"; Value = "-1.41%" }
    @{ Cell = "D15"; Value = "0.0006035" }
    @{ Cell = "E15"; Value = "-0.28%" }
    @{ Cell = "D16"; Value = "0.006197" }
    @{ Cell = "E16"; Value = "1.82%" }
    @{ Cell = "D17"; Value = "3.448" }
    @{ Cell = "E17"; Value = "-1.91%" }
    @{ Cell = "D18"; Value = "3.149" }
    @{ Cell = "E18"; Value = "-1.47%" }
    @{ Cell = "D19"; Value = "2.180" }
    @{ Cell = "E19"; Value = "-1.02%" }
    @{ Cell = "E20"; Value = "1.65%" }
    @{ Cell = "E21"; Value = "0.86%" }
    @{ Cell = "D22"; Value = "4.038" }
    @{ Cell = "E22"; Value = "4.97%" }
    @{ Cell = "D23"; Value = "0.04234" }
    @{ Cell = "E23"; Value = "0.67%" }
    @{ Cell = "E24"; Value = "-0.46%" }
    @{ Cell = "D25"; Value = "0.004024" }
    @{ Cell = "E25"; Value = "-5.93%" }
    @{ Cell = "E26"; Value = "-0.04%" }
    @{ Cell = "E27"; Value = "-11.66%" }
    @{ Cell = "D40"; Value = "0.03827" }
    @{ Cell = "E40"; Value = "-0.04%" }
    @{ Cell = "E41"; Value = "1.18%" }
    @{ Cell = "D42"; Value = "0.003767" }
    @{ Cell = "E42"; Value = "-39.56%" }
    @{ Cell = "D43"; Value = "0.002430" }
    @{ Cell = "E43"; Value = "10.41%" }
    @{ Cell = "D44"; Value = "0.01517" }
    @{ Cell = "E44"; Value = "33.85%" }
    @{ Cell = "D45"; Value = "0.00005100" }
    @{ Cell = "E45"; Value = "-6.24%" }
    @{ Cell = "E46"; Value = "-0.06%" }
    @{ Cell = "E47"; Value = "-17.05%" }
    @{ Cell = "D48"; Value = "0.1165" }
    @{ Cell = "E48"; Value = "19.70%" }
    @{ Cell = "E49"; Value = "-0.06%" }
    @{ Cell = "E50"; Value = "-0.06%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
